# Updated cryptos list on Tue Apr 11 07:44:15 UTC 2023 with GitHub Actions
# The "Price" column (D) holds values that look numeric (e.g. "1.0000",
# "43.47") but must stay plain text, exactly like the source workbook.
# Prefixing with a leading apostrophe forces Excel to store them as text
# (quote-prefixed) instead of silently coercing them into numbers, which
# would corrupt the fidelity of the values (dropping trailing zeros etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "'30.103.90"
$ws.Range("E2").Value = "  +5.70%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "'1.923.24"
$ws.Range("E3").Value = "  +2.75%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.80%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'331.92"
$ws.Range("E5").Value = "  +5.12%  "

# --- Row 6: USDC ---
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.79%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "'0.5244"
$ws.Range("E7").Value = "  +2.99%  "

# --- Row 8: Cardano ---
$ws.Range("D8").Value = "'0.4076"
$ws.Range("E8").Value = "  +4.62%  "

# --- Row 9: Dogecoin ---
$ws.Range("D9").Value = "'0.08554"
$ws.Range("E9").Value = "  +2.57%  "

# --- Row 10: OKB ---
$ws.Range("D10").Value = "'43.47"

# --- Row 11: Polygon ---
$ws.Range("E11").Value = "  +2.54%  "

# --- Row 12: Solana ---
$ws.Range("E12").Value = "  +10.33%  "

# --- Row 13: Polkadot ---
$ws.Range("D13").Value = "'6.434"
$ws.Range("E13").Value = "  +3.34%  "

# --- Row 14: WrappedEther ---
$ws.Range("D14").Value = "'1.919.95"
$ws.Range("E14").Value = "  +2.63%  "

# --- Row 15: Chainlink ---
$ws.Range("D15").Value = "'7.415"
$ws.Range("E15").Value = "  +2.04%  "

# --- Row 16: BinanceUSD ---
$ws.Range("E16").Value = "  -0.83%  "

# --- Row 17: Litecoin ---
$ws.Range("D17").Value = "'96.96"
$ws.Range("E17").Value = "  +6.37%  "

# --- Row 18: ShibaInu ---
$ws.Range("E18").Value = "  +1.20%  "

# --- Row 19: TRON ---
$ws.Range("D19").Value = "'0.06718"
$ws.Range("E19").Value = "  -0.14%  "

# --- Row 20: Avalanche ---
$ws.Range("E20").Value = "  +3.73%  "

# --- Row 21: Dai ---
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.67%  "

# --- Row 22: Uniswap ---
$ws.Range("D22").Value = "'6.069"
$ws.Range("E22").Value = "  +2.74%  "

# --- Row 23: WrappedBTC ---
$ws.Range("D23").Value = "'30.121.12"
$ws.Range("E23").Value = "  +5.67%  "

# --- Row 24: Cosmos ---
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "  +1.68%  "

# --- Row 25: Toncoin ---
$ws.Range("D25").Value = "'2.224"
$ws.Range("E25").Value = "  -0.06%  "

# --- Row 26: WrappedliquidstakedEther2.0 ---
$ws.Range("D26").Value = "'2.141.09"
$ws.Range("E26").Value = "  +2.75%  "

# --- Row 27: EthereumClassic ---
$ws.Range("D27").Value = "'21.19"
$ws.Range("E27").Value = "  +2.76%  "

# --- Row 28: Monero ---
$ws.Range("D28").Value = "'160.23"
$ws.Range("E28").Value = "  -1.04%  "

# --- Row 29: LidoDAOToken ---
$ws.Range("D29").Value = "'2.478"
$ws.Range("E29").Value = "  +3.07%  "

# --- Row 30: BitcoinCash ---
$ws.Range("D30").Value = "'129.67"
$ws.Range("E30").Value = "  +2.77%  "

# --- Row 31: ImmutableX ---
$ws.Range("E31").Value = "  +4.85%  "

# --- Row 32: Stellar ---
$ws.Range("D32").Value = "'0.1058"
$ws.Range("E32").Value = "  +1.54%  "

# --- Row 33: Filecoin ---
$ws.Range("D33").Value = "'6.127"
$ws.Range("E33").Value = "  +6.67%  "

# --- Row 34: HuobiToken ---
$ws.Range("D34").Value = "'3.649"
$ws.Range("E34").Value = "  +1.08%  "

# --- Row 35: VeChain ---
$ws.Range("D35").Value = "'0.02523"
$ws.Range("E35").Value = "  +2.90%  "

# --- Row 36: Hedera ---
$ws.Range("D36").Value = "'0.06621"
$ws.Range("E36").Value = "  +1.19%  "

# --- Row 37: Algorand ---
$ws.Range("D37").Value = "'0.2228"
$ws.Range("E37").Value = "  +3.35%  "

# --- Rows 38 & 39: FraxShare and ARBITRUM swap places ---
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.241"
$ws.Range("E38").Value = "  +4.82%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'9.052"
$ws.Range("E39").Value = "  +2.60%  "

# --- Row 41: TheSandbox ---
$ws.Range("D41").Value = "'0.6576"
$ws.Range("E41").Value = "  +3.28%  "

# --- Row 42: Aptos ---
$ws.Range("D42").Value = "'11.70"
$ws.Range("E42").Value = "  +5.70%  "

# --- Row 43: TrustWalletToken ---
$ws.Range("D43").Value = "'1.241"
$ws.Range("E43").Value = "  -0.25%  "

# --- Row 44: Decentraland ---
$ws.Range("D44").Value = "'0.6205"
$ws.Range("E44").Value = "  +3.64%  "

# --- Row 45: EnergySwap ---
$ws.Range("D45").Value = "'13.30"
$ws.Range("E45").Value = "  +1.96%  "

# --- Row 46: PancakeSwap ---
$ws.Range("D46").Value = "'3.787"
$ws.Range("E46").Value = "  +2.66%  "

# --- Row 47: NEARProtocol ---
$ws.Range("D47").Value = "'2.094"
$ws.Range("E47").Value = "  +4.69%  "

# --- Row 48: EOS ---
$ws.Range("E48").Value = "  +2.93%  "

# --- Row 49: Quant ---
$ws.Range("D49").Value = "'125.23"
$ws.Range("E49").Value = "  +2.83%  "

# --- Rows 50 & 51: Aave and WEMIXTOKEN swap places ---
$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'1.166"
$ws.Range("E50").Value = "  +1.42%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'80.00"
$ws.Range("E51").Value = "  +5.11%  "
